$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '25.805.32'
$c.ClearFormats()
$ws.Range('E2').Value = '  +0.48%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.750.68'
$c.ClearFormats()
$ws.Range('E3').Value = '  +0.36%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.ClearFormats()
$ws.Range('E4').Value = '  -0.02%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '236.18'
$c.ClearFormats()
$ws.Range('E6').Value = '  -0.04%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.5084'
$c.ClearFormats()
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.2688'
$c.ClearFormats()
$ws.Range('E8').Value = '  +7.61%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.06196'
$c.ClearFormats()
$ws.Range('E9').Value = '  +4.40%  '
$ws.Range('B10').Value = 'WrappedEther'
$ws.Range('C10').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '1.760.42'
$c.ClearFormats()
$ws.Range('E10').Value = '  +0.91%  '
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.06934'
$c.ClearFormats()
$ws.Range('E11').Value = '  +2.31%  '
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '15.44'
$c.ClearFormats()
$ws.Range('E12').Value = '  +4.76%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.6283'
$c.ClearFormats()
$ws.Range('E13').Value = '  +10.68%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '4.480'
$c.ClearFormats()
$ws.Range('E14').Value = '  +0.56%  '
$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '77.94'
$c.ClearFormats()
$ws.Range('E15').Value = '  +0.96%  '
$ws.Range('B16').Value = 'BinanceUSD'
$ws.Range('C16').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.ClearFormats()
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('B17').Value = 'Dai'
$ws.Range('C17').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('E17').Value = '  -0.07%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '25.822.00'
$c.ClearFormats()
$ws.Range('E18').Value = '  +0.31%  '
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '11.66'
$c.ClearFormats()
$ws.Range('E19').Value = '  +1.59%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '0.000006694'
$c.ClearFormats()
$ws.Range('E20').Value = '  +2.08%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '1.977.33'
$c.ClearFormats()
$ws.Range('E21').Value = '  +0.63%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '4.061'
$c.ClearFormats()
$ws.Range('E22').Value = '  +2.42%  '
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '8.253'
$c.ClearFormats()
$ws.Range('E23').Value = '  +4.96%  '
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '5.139'
$c.ClearFormats()
$ws.Range('E24').Value = '  +2.50%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '136.70'
$c.ClearFormats()
$ws.Range('E25').Value = '  +0.51%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '1.453'
$c.ClearFormats()
$ws.Range('E26').Value = '  -1.49%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '15.14'
$c.ClearFormats()
$ws.Range('E27').Value = '  +3.48%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '1.748'
$c.ClearFormats()
$ws.Range('E28').Value = '  -4.03%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '102.60'
$c.ClearFormats()
$ws.Range('E29').Value = '  +0.70%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '0.08194'
$c.ClearFormats()
$ws.Range('E30').Value = '  +1.66%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '3.709'
$c.ClearFormats()
$ws.Range('E31').Value = '  -1.44%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '3.406'
$c.ClearFormats()
$ws.Range('E32').Value = '  +2.65%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.04442'
$c.ClearFormats()
$ws.Range('E33').Value = '  +1.12%  '
$ws.Range('B34').Value = 'Frax'
$ws.Range('C34').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '1.000'
$c.ClearFormats()
$ws.Range('E34').Value = '  -0.07%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '2.648'
$c.ClearFormats()
$ws.Range('E35').Value = '  +1.56%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.9995'
$c.ClearFormats()
$ws.Range('E36').Value = '  +2.67%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.6024'
$c.ClearFormats()
$ws.Range('E37').Value = '  +0.52%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '2.646'
$c.ClearFormats()
$ws.Range('E38').Value = '  -1.29%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.01567'
$c.ClearFormats()
$ws.Range('E39').Value = '  +4.90%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '1.947'
$c.ClearFormats()
$ws.Range('E40').Value = '  -3.85%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.ClearFormats()
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('B42').Value = 'PaxosStandard'
$ws.Range('C42').Value = 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '1.002'
$c.ClearFormats()
$ws.Range('E42').Value = '  +0.00%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '101.26'
$c.ClearFormats()
$ws.Range('E43').Value = '  -2.17%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.3835'
$c.ClearFormats()
$ws.Range('E44').Value = '  +3.32%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.7509'
$c.ClearFormats()
$ws.Range('E45').Value = '  -0.68%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '4.900'
$c.ClearFormats()
$ws.Range('E46').Value = '  -5.04%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.05507'
$c.ClearFormats()
$ws.Range('E47').Value = '  +7.74%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.1102'
$c.ClearFormats()
$ws.Range('E48').Value = '  +3.21%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '5.961'
$c.ClearFormats()
$ws.Range('E49').Value = '  +1.65%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '30.14'
$c.ClearFormats()
$ws.Range('E50').Value = '  +0.04%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '52.80'
$c.ClearFormats()
$ws.Range('E51').Value = '  +0.66%  '
